# Generate Report for Handback
# Update the "generate date" / handoff / handback timestamp values.
# These cells are stored as plain text strings (formatted to look like
# dates via a custom number format), so we set them explicitly as text
# to avoid Excel re-interpreting/reformatting them as real date values.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-31 21:20:51"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-31 21:20:45"
$wsZhCn.Range("K2").Value = "2016-08-31 21:21:08"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-31 21:21:18"
